# Apply test-row fill-in edits to the "Test Case Template 1" sheet.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Test Case Template 1")

# For rows 14-20: set the Actual Result (D), Pass/Fail (E), Browser (G)
# and Device (H) columns to the same values used across the test run.
for ($r = 14; $r -le 20; $r++) {
    $ws.Cells.Item($r, 4).Value = "As Expected"   # D: Actual Result
    $ws.Cells.Item($r, 5).Value = "Pass"          # E: Pass/Fail
    $ws.Cells.Item($r, 7).Value = "Chrome"        # G: Browser Tested on
    $ws.Cells.Item($r, 8).Value = "Desktop"       # H: Device tested on
}

# The "Browser Tested on" (G) and "Device tested on" (H) list validations
# originally only covered rows 14:16; extend them down to the newly
# filled-in rows (14:20) to match the dropdown used for data entry above.
$ws.Range("G14:G20").Validation.Delete() | Out-Null
$ws.Range("G14:G20").Validation.Add(3, 1, 1, "=Sheet1!`$B`$1:`$B`$11") | Out-Null
$ws.Range("H14:H20").Validation.Delete() | Out-Null
$ws.Range("H14:H20").Validation.Add(3, 1, 1, "=Sheet1!`$C`$1:`$C`$5") | Out-Null

# Move the active selection to reflect where the author ended up.
$ws.Range("F18").Select() | Out-Null
